$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2503.481290811585, 1699.8523375677767, 1673.9643730521309),
    @(2277.9971094061325, 1501.9362601401947, 1402.0288730485797),
    @(2560.1242476724956, 1737.303560732443, 1580.4586850596445),
    @(2494.2668094583914, 1876.6353147064599, 1893.508775158972),
    @(2681.663715749919, 1754.5415264857279, 1661.4949790064081),
    @(2478.7615939391708, 1810.0447733893373, 1876.8746931705411),
    @(2396.3740662144251, 1851.2021005594868, 1664.0425653320099),
    @(2466.9899163351297, 1937.8683342009456, 1847.6799461779342),
    @(2643.8902174177633, 1944.4321728919576, 1636.5747422631507),
    @(2394.9207777787101, 1496.7314448948889, 1449.4932515435044),
    @(2163.7693144377299, 1575.0359845990301, 1427.7475679716438),
    @(2859.2191420186659, 2218.1049445540334, 1914.8994363674237),
    @(2528.1905879101214, 1932.435834573078, 1794.9433085700407),
    @(2610.7289789880551, 2074.3284564691662, 2050.6666685628061),
    @(2531.3176201022948, 1976.9826572008753, 1827.4805076693092),
    @(2592.6392331263787, 1797.0722032246758, 1550.1668649135759)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
